$d = $word.ActiveDocument

$replacements = @(
    @{old="97×61="; new="63×69="},
    @{old="79×42="; new="98×32="},
    @{old="35×27="; new="24×94="},
    @{old="97×29="; new="59×97="},
    @{old="80×33="; new="63×84="},
    @{old="94×63="; new="34×57="},
    @{old="70×25="; new="53×69="},
    @{old="57×72="; new="33×78="},
    @{old="34×90="; new="31×95="},
    @{old="96×68="; new="34×18="},
    @{old="73×54="; new="53×43="},
    @{old="22×51="; new="82×65="},
    @{old="51×42="; new="49×78="},
    @{old="62×81="; new="93×58="},
    @{old="83×30="; new="87×93="},
    @{old="39×46="; new="49×26="},
    @{old="92×23="; new="48×32="},
    @{old="68×79="; new="44×68="},
    @{old="88×54="; new="73×95="},
    @{old="26×31="; new="35×28="},
    @{old="64×18="; new="95×92="},
    @{old="33×65="; new="30×43="},
    @{old="53×13="; new="41×16="},
    @{old="80×25="; new="71×75="},
    @{old="86×99="; new="17×11="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
